$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices) are written as literal text,
# matching the source data which stores these as strings (e.g. "70.517.18").

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.517.18'
$ws.Range('E2').Value = '  +0.78%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.626.24'
$ws.Range('E3').Value = '  +2.42%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.80'
$ws.Range('E5').Value = '  +0.07%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '196.46'
$ws.Range('E6').Value = '  +0.05%  '

# Row 7
$ws.Range('E7').Value = '  -0.61%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.02%  '

# Row 9
$ws.Range('E9').Value = '  +0.37%  '

# Row 10
$ws.Range('E10').Value = '  -0.96%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.35'
$ws.Range('E11').Value = '  -1.04%  '

# Row 12
$ws.Range('E12').Value = '  -0.09%  '

# Row 13
$ws.Range('E13').Value = '  +0.49%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.205.72'
$ws.Range('E14').Value = '  +2.50%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '601.24'
$ws.Range('E15').Value = '  -1.38%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.98'
$ws.Range('E16').Value = '  +0.77%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.637.37'
$ws.Range('E17').Value = '  +0.68%  '

# Row 18
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.639.17'
$ws.Range('E18').Value = '  +2.51%  '

# Row 19
$ws.Range('E19').Value = '  -0.81%  '

# Row 20
$ws.Range('E20').Value = '  +1.67%  '

# Row 21
$ws.Range('E21').Value = '  +0.14%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.21'
$ws.Range('E22').Value = '  +1.26%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.19'
$ws.Range('E23').Value = '  -2.26%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '102.15'
$ws.Range('E24').Value = '  -0.41%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.62'
$ws.Range('E25').Value = '  -0.11%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.00'
$ws.Range('E26').Value = '  -3.66%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.67'
$ws.Range('E27').Value = '  -2.28%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.68'
$ws.Range('E28').Value = '  +0.65%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.85'
$ws.Range('E29').Value = '  +0.86%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.65'
$ws.Range('E30').Value = '  +7.21%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.29'
$ws.Range('E31').Value = '  +2.35%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.28'
$ws.Range('E32').Value = '  -2.92%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.118'
$ws.Range('E33').Value = '  +2.30%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.46'
$ws.Range('E34').Value = '  +0.33%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0889'
$ws.Range('E35').Value = '  +3.63%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.928.93'
$ws.Range('E36').Value = '  +4.97%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '539.62'
$ws.Range('E37').Value = '  +9.17%  '

# Row 38
$ws.Range('E38').Value = '  +0.08%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.07'
$ws.Range('E39').Value = '  +0.65%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.94'
$ws.Range('E40').Value = '  +0.74%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.389'
$ws.Range('E41').Value = '  -1.24%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.55'
$ws.Range('E42').Value = '  -2.77%  '

# Row 43
$ws.Range('E43').Value = '  +0.15%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0462'
$ws.Range('E44').Value = '  +1.21%  '

# Row 45
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.40'
$ws.Range('E45').Value = '  +2.69%  '

# Row 46
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.88'
$ws.Range('E46').Value = '  +1.26%  '

# Row 47
$ws.Range('E47').Value = '  -0.23%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.58'
$ws.Range('E48').Value = '  -0.79%  '

# Row 49
$ws.Range('E49').Value = '  -0.25%  '

# Row 50
$ws.Range('E50').Value = '  -2.05%  '

# Row 51
$ws.Range('E51').Value = '  +1.67%  '
